# Append the new daily data row (2025/10/01) to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 42

# Column A holds dates stored as plain text (e.g. "2025/09/30"), not real
# Excel dates. Force text formatting before assigning the value so Excel's
# COM layer doesn't auto-convert the "2025/10/01" string into a date serial
# number, then restore the default "Normal" style so the cell doesn't pick
# up a stray number-format style that the other rows don't have.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/10/01"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "水"
$ws.Cells.Item($newRow, 3).Value = 0
$ws.Cells.Item($newRow, 4).Value = 11
